$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D values look numeric (e.g. "59.01"); Excel would auto-convert them to
# numbers on assignment, losing the intended text representation. Force the
# "Price" column cells being updated to a text number format first, matching
# the original workbook where these values are stored as text.
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D16","D17","D18","D19","D20","D21","D22","D23","D24","D26","D27","D28","D29","D32","D35","D38","D42","D44","D47","D48","D49","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '36.475.64'
$ws.Range("E2").Value = '  -2.72%  '

$ws.Range("D3").Value = '1.984.64'
$ws.Range("E3").Value = '  -3.41%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '244.06'
$ws.Range("E5").Value = '  -3.37%  '

$ws.Range("D6").Value = '0.628'
$ws.Range("E6").Value = '  -4.31%  '

$ws.Range("D7").Value = '59.01'
$ws.Range("E7").Value = '  -10.31%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '0.375'
$ws.Range("E9").Value = '  -2.02%  '

$ws.Range("D10").Value = '57.64'
$ws.Range("E10").Value = '  -3.60%  '

$ws.Range("D11").Value = '0.0823'
$ws.Range("E11").Value = '  +7.02%  '

$ws.Range("E12").Value = '  -0.98%  '

$ws.Range("D13").Value = '23.72'
$ws.Range("E13").Value = '  +5.32%  '

$ws.Range("D14").Value = '0.864'
$ws.Range("E14").Value = '  -6.15%  '

$ws.Range("E15").Value = '  -5.65%  '

$ws.Range("D16").Value = '2.275.89'
$ws.Range("E16").Value = '  -3.31%  '

$ws.Range("D17").Value = '5.46'
$ws.Range("E17").Value = '  -2.41%  '

$ws.Range("D18").Value = '1.981.68'
$ws.Range("E18").Value = '  -3.41%  '

$ws.Range("D19").Value = '36.339.02'
$ws.Range("E19").Value = '  -2.65%  '

$ws.Range("D20").Value = '70.55'
$ws.Range("E20").Value = '  -4.40%  '

$ws.Range("D21").Value = '0.0₃0863'
$ws.Range("E21").Value = '  -1.47%  '

$ws.Range("D22").Value = '5.33'
$ws.Range("E22").Value = '  -2.67%  '

$ws.Range("D23").Value = '233.60'
$ws.Range("E23").Value = '  -2.87%  '

$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("E25").Value = '  -1.15%  '

$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  -4.10%  '

$ws.Range("D27").Value = '10.08'
$ws.Range("E27").Value = '  +0.20%  '

$ws.Range("D28").Value = '162.12'
$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("D29").Value = '19.84'
$ws.Range("E29").Value = '  -0.91%  '

$ws.Range("E30").Value = '  +11.06%  '

$ws.Range("E31").Value = '  -2.02%  '

$ws.Range("D32").Value = '1.20'
$ws.Range("E32").Value = '  -0.96%  '

$ws.Range("E33").Value = '  -7.11%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("D35").Value = '4.42'
$ws.Range("E35").Value = '  -6.08%  '

$ws.Range("E36").Value = '  +4.41%  '

$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").Value = '2.27'
$ws.Range("E38").Value = '  -7.58%  '

$ws.Range("E39").Value = '  -3.49%  '

$ws.Range("E40").Value = '  +2.01%  '

$ws.Range("E41").Value = '  -0.71%  '

$ws.Range("D42").Value = '0.0959'
$ws.Range("E42").Value = '  -7.42%  '

$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("D44").Value = '0.0213'
$ws.Range("E44").Value = '  -2.49%  '

$ws.Range("E45").Value = '  -5.13%  '

$ws.Range("E46").Value = '  -4.37%  '

$ws.Range("D47").Value = '16.19'
$ws.Range("E47").Value = '  -5.80%  '

$ws.Range("D48").Value = '1.383.66'
$ws.Range("E48").Value = '  -2.70%  '

$ws.Range("D49").Value = '7.50'
$ws.Range("E49").Value = '  -6.12%  '

$ws.Range("D50").Value = '2.85'
$ws.Range("E50").Value = '  -3.35%  '

$ws.Range("D51").Value = '45.54'
$ws.Range("E51").Value = '  -2.98%  '
